# 4.2 eInk switched to gray2 mode
#
# Typography sheet: bump Bpp from 1 to 2 for the three fonts in use,
# and adjust their pixel Size to match the new gray2 rendering.
#
# Translation sheet: rework two strings ("Upd" alignment/GB text and the
# "Temperature" label), and add a new "SingleUseId27" text row that keeps
# the original "read" wording.

$wb = $excel.ActiveWorkbook

$typography = $wb.Worksheets.Item("Typography")

$typography.Range("D4").Value = 24
$typography.Range("E4").Value = 2

$typography.Range("D5").Value = 52
$typography.Range("E5").Value = 2

$typography.Range("D6").Value = 80
$typography.Range("E6").Value = 2

$translation = $wb.Worksheets.Item("Translation")

$translation.Range("D5").Value = "Center"
$translation.Range("F5").Value = "refresh"

$translation.Range("F19").Value = "Temp."

$translation.Range("B23").Value = "SingleUseId27"
$translation.Range("C23").Value = "Default"
$translation.Range("D23").Value = "Left"
$translation.Range("E23").Value = "LTR"
$translation.Range("F23").Value = "read"
